$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply "0.00" number format to G34:G36 (dedup'd automatically) ---
$ws.Range("G34").NumberFormat = "0.00"
$ws.Range("G35").NumberFormat = "0.00"
$ws.Range("G36").NumberFormat = "0.00"

# --- 2. Add the two new data rows (52 & 53) ---
$ws.Range("B52").Value = "IL"
$ws.Range("C52").Value = "milton township"
$ws.Range("D52").Value = 1749451
$ws.Range("E52").Value = 120211
$ws.Range("F52").Value = 41.858611000000003
$ws.Range("G52").Value = -88.089444

$ws.Range("B53").Value = "IL"
$ws.Range("C53").Value = "lisle township"
$ws.Range("D53").Value = 1743952
$ws.Range("E53").Value = 119045
$ws.Range("F53").Value = 41.771667000000001
$ws.Range("G53").Value = -88.088611

# --- 3. Row heights for the two new rows ---
$ws.Rows.Item(52).RowHeight = 20
$ws.Rows.Item(53).RowHeight = 20

# --- 4. Fonts: columns B,C,D,E,G get Lucida Grande 11; column F gets Arial 16 FF001D35 ---
$rngSmall = $ws.Range("B52:E53")
$rngSmall.Font.Size = 11
$rngSmall.Font.Name = "Lucida Grande"

$rngG = $ws.Range("G52:G53")
$rngG.Font.Size = 11
$rngG.Font.Name = "Lucida Grande"

$rngF = $ws.Range("F52:F53")
$rngF.Font.Size = 16
$rngF.Font.Name = "Arial"
$rngF.Font.Color = 3480832

# --- 5. Column widths ---
$ws.Columns.Item(3).ColumnWidth = 58.6
$ws.Columns.Item(7).ColumnWidth = 71.5

# --- 6. Sheet view / selection ---
$ws.Range("F57").Select() | Out-Null

Write-Host "edit complete"
